$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.458.74'
$ws.Range("E2").Value = '  -4.21%  '
$ws.Range("D3").Value = '2.374.84'
$ws.Range("E3").Value = '  -5.16%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '502.29'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -6.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.95%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.552'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.99%  '
$ws.Range("D9").Value = '2.395.28'
$ws.Range("E9").Value = '  -4.49%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0956'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.00%  '
$ws.Range("E11").Value = '  -1.62%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.317'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.80%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.60'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -11.23%  '
$ws.Range("D14").Value = '2.797.60'
$ws.Range("E14").Value = '  -5.10%  '
$ws.Range("D15").Value = '56.968.45'
$ws.Range("E15").Value = '  -3.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.52'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.88%  '
$ws.Range("E17").Value = '  -3.58%  '
$ws.Range("D18").Value = '2.379.78'
$ws.Range("E18").Value = '  -4.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.12'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '310.50'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.33%  '
$ws.Range("E21").Value = '  -5.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.20'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.26%  '
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.51'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.53%  '
$ws.Range("E25").Value = '  +0.69%  '
$ws.Range("D26").Value = '2.497.94'
$ws.Range("E26").Value = '  -4.64%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.371'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -9.48%  '
$ws.Range("E28").Value = '  -6.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.20'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.56%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '174.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.66'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.38%  '
$ws.Range("D32").Value = '0.0₃0709'
$ws.Range("E32").Value = '  -6.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.10'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.99%  '
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("E35").Value = '  -7.70%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.994'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.33%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.74'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.10%  '
$ws.Range("E38").Value = '  -1.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.74'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.67%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '35.80'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.42'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.769'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '129.42'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.29%  '
$ws.Range("E44").Value = '  -3.88%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.78'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.95%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.571'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '253.66'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.76%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0897'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0483'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.52%  '
$ws.Range("E50").Value = '  -4.98%  '
$ws.Range("E51").Value = '  -5.36%  '
